$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header row (row 1), pushing existing
# data down. This mirrors inserting "Presidencia 06" and "Presidencia 12"
# at the top of the data table.
$ws.Rows("2:3").Insert()

$ws.Range("A3").Value = "Presidencia 12"
$ws.Range("B3").Value = "pr_12"
$ws.Range("C3").Value = "#d6ccc2"

$ws.Range("A2").Value = "Presidencia 06"
$ws.Range("B2").Value = "pr_06"
$ws.Range("C2").Value = "#b5838d"

$ws.Range("C2").Select()
